$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2039.5
$ws.Range("J17").Value = 2039.5
$ws.Range("L17").Value = 6118.5
$ws.Range("N17").Value = -6454.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2766.6667
$ws.Range("I28").Value = 2120
$ws.Range("K28").Value = 2120
$ws.Range("M28").Value = -1635

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 11364930
$ws.Range("I80").Value = 20834326
$ws.Range("K80").Value = 62502978
$ws.Range("M80").Value = -62501980

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 11364930
$ws.Range("I83").Value = 20834326
$ws.Range("K83").Value = 187508934
$ws.Range("M83").Value = -187503942

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 5210.4287
$ws.Range("I96").Value = 829.1818
$ws.Range("K96").Value = 2487.5454
$ws.Range("M96").Value = -1114.5454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3390.8
$ws.Range("I106").Value = 3156.6
$ws.Range("K106").Value = 3156.6
$ws.Range("M106").Value = -2525.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1158547.9
$ws.Range("I116").Value = 1990688.8
$ws.Range("J116").Value = 6352.846
$ws.Range("K116").Value = 1990688.8
$ws.Range("L116").Value = 6352.846
$ws.Range("M116").Value = -1987246.8
$ws.Range("N116").Value = -13236.846

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4866.3213
$ws.Range("J138").Value = 5787.8887
$ws.Range("L138").Value = 17363.6661
$ws.Range("N138").Value = -27643.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17974.012
$ws.Range("I32").Value = 16207.108
$ws.Range("K32").Value = 16207.108
$ws.Range("M32").Value = -15920.108

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 32771.5
$ws.Range("J37").Value = 48600
$ws.Range("L37").Value = 48600
$ws.Range("N37").Value = -49146

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 50000
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50630

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 63499.5
$ws.Range("J80").Value = 67999.336
$ws.Range("L80").Value = 67999.336
$ws.Range("N80").Value = -69995.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 63499.5
$ws.Range("J83").Value = 67999.336
$ws.Range("L83").Value = 203998.008
$ws.Range("N83").Value = -213982.008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1543.5714
$ws.Range("I110").Value = 1633.75
$ws.Range("K110").Value = 1633.75
$ws.Range("M110").Value = 411.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 8138
$ws.Range("I122").Value = 6810.75
$ws.Range("K122").Value = 20432.25
$ws.Range("M122").Value = -17982.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 13963.304
$ws.Range("I132").Value = 17409.129
$ws.Range("K132").Value = 52227.387
$ws.Range("M132").Value = -49697.387

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5186.878
$ws.Range("I20").Value = 4889.1304
$ws.Range("J20").Value = 5567.3335
$ws.Range("K20").Value = 4889.1304
$ws.Range("L20").Value = 5567.3335
$ws.Range("M20").Value = -4642.1304
$ws.Range("N20").Value = -6061.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 23793.572
$ws.Range("J82").Value = 70400
$ws.Range("L82").Value = 70400
$ws.Range("N82").Value = -71166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 23793.572
$ws.Range("J85").Value = 70400
$ws.Range("L85").Value = 70400
$ws.Range("N85").Value = -73052

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6813.2856
$ws.Range("I86").Value = 4349.5
$ws.Range("J86").Value = 7798.8
$ws.Range("K86").Value = 4349.5
$ws.Range("L86").Value = 7798.8
$ws.Range("M86").Value = -3226.5
$ws.Range("N86").Value = -10044.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 6813.2856
$ws.Range("I89").Value = 4349.5
$ws.Range("J89").Value = 7798.8
$ws.Range("K89").Value = 21747.5
$ws.Range("L89").Value = 38994
$ws.Range("M89").Value = -16131.5
$ws.Range("N89").Value = -50226

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3879.6
$ws.Range("I107").Value = 1349.5
$ws.Range("J107").Value = 14000
$ws.Range("K107").Value = 1349.5
$ws.Range("L107").Value = 14000
$ws.Range("M107").Value = 570.5
$ws.Range("N107").Value = -17840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 100000000
$ws.Range("I3").Value = 100000000
$ws.Range("K3").Value = 100000000
$ws.Range("M3").Value = -99999887

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23258964
$ws.Range("I31").Value = 29414068
$ws.Range("K31").Value = 29414068
$ws.Range("M31").Value = -29413773

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 23258964
$ws.Range("I34").Value = 29414068
$ws.Range("K34").Value = 29414068
$ws.Range("M34").Value = -29413866

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1361.2307
$ws.Range("I94").Value = 350.1111
$ws.Range("K94").Value = 350.1111
$ws.Range("M94").Value = 100.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 6697.1113
$ws.Range("I99").Value = 3779.8
$ws.Range("K99").Value = 3779.8
$ws.Range("M99").Value = -2281.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1130.3529
$ws.Range("I107").Value = 865.7857
$ws.Range("K107").Value = 865.7857
$ws.Range("M107").Value = 1054.2143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 6697.1113
$ws.Range("I126").Value = 3779.8
$ws.Range("K126").Value = 11339.4
$ws.Range("M126").Value = -8869.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 78743.60000000001
$ws.Range("J133").Value = 78743.60000000001
$ws.Range("L133").Value = 78743.60000000001
$ws.Range("N133").Value = -83803.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 296297.38
$ws.Range("J141").Value = 307844.3
$ws.Range("L141").Value = 307844.3
$ws.Range("N141").Value = -318204.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 653
$ws.Range("I13").Value = 72.666664
$ws.Range("J13").Value = 1233.3334
$ws.Range("K13").Value = 217.999992
$ws.Range("L13").Value = 3700.0002
$ws.Range("M13").Value = -49.99999199999999
$ws.Range("N13").Value = -4036.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 579.4286
$ws.Range("I107").Value = 369.75
$ws.Range("J107").Value = 628.7646999999999
$ws.Range("K107").Value = 1109.25
$ws.Range("L107").Value = 1886.2941
$ws.Range("M107").Value = 810.75
$ws.Range("N107").Value = -5726.2941

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2016.5358
$ws.Range("I113").Value = 617
$ws.Range("J113").Value = 2249.7917
$ws.Range("K113").Value = 1851
$ws.Range("L113").Value = 6749.375100000001
$ws.Range("M113").Value = 319
$ws.Range("N113").Value = -11089.3751

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 20000
$ws.Range("I49").Value = 20000
$ws.Range("J49").Value = 20000
$ws.Range("K49").Value = 20000
$ws.Range("L49").Value = 20000
$ws.Range("M49").Value = -19816
$ws.Range("N49").Value = -20368

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 34011564
$ws.Range("I102").Value = 72873640
$ws.Range("K102").Value = 72873640
$ws.Range("M102").Value = -72872018

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1622
$ws.Range("I113").Value = 1622
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1622
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 548
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4149.054
$ws.Range("I122").Value = 2075.1765
$ws.Range("K122").Value = 6225.529500000001
$ws.Range("M122").Value = -3775.529500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 37269.777
$ws.Range("I132").Value = 44877.914
$ws.Range("J132").Value = 4762.273
$ws.Range("K132").Value = 134633.742
$ws.Range("L132").Value = 14286.819
$ws.Range("M132").Value = -132103.742
$ws.Range("N132").Value = -19346.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5010.926
$ws.Range("I16").Value = 4608.909
$ws.Range("K16").Value = 4608.909
$ws.Range("M16").Value = -4438.909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 36464760
$ws.Range("I40").Value = 20838358
$ws.Range("J40").Value = 83343960
$ws.Range("K40").Value = 20838358
$ws.Range("L40").Value = 83343960
$ws.Range("M40").Value = -20838222
$ws.Range("N40").Value = -83344232

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5394
$ws.Range("I46").Value = 5873.2354
$ws.Range("J46").Value = 4767.3076
$ws.Range("K46").Value = 5873.2354
$ws.Range("L46").Value = 4767.3076
$ws.Range("M46").Value = -5685.2354
$ws.Range("N46").Value = -5143.3076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2477.7368
$ws.Range("J93").Value = 3289.1
$ws.Range("L93").Value = 3289.1
$ws.Range("N93").Value = -5785.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 78000
$ws.Range("J121").Value = 78000
$ws.Range("L121").Value = 78000
$ws.Range("N121").Value = -81494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 5294.75
$ws.Range("J130").Value = 5294.75
$ws.Range("L130").Value = 5294.75
$ws.Range("N130").Value = -15334.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4734.21
$ws.Range("I132").Value = 4073.8975
$ws.Range("K132").Value = 12221.6925
$ws.Range("M132").Value = -9691.692500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6126.857
$ws.Range("I136").Value = 5997
$ws.Range("J136").Value = 6224.25
$ws.Range("K136").Value = 17991
$ws.Range("L136").Value = 18672.75
$ws.Range("M136").Value = -15441
$ws.Range("N136").Value = -23772.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 40497
$ws.Range("I37").Value = 34990
$ws.Range("K37").Value = 34990
$ws.Range("M37").Value = -34787

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3746.8293
$ws.Range("I132").Value = 3592.4707
$ws.Range("J132").Value = 3856.1667
$ws.Range("K132").Value = 10777.4121
$ws.Range("L132").Value = 11568.5001
$ws.Range("M132").Value = -8247.4121
$ws.Range("N132").Value = -16628.5001
